# TeamANSReport.xlsx - "Change to Excel for sprint 2 stories"
#
# Adds the two Sprint-2 stories (US01 "Dates before current date" and
# US04 "Marriage before divorce", both owned by SK) to the Backlog sheet
# and to the Sprint2 detail sheet, and updates the Burndown sheet's LOC
# figure for the latest check-in.

$wb = $excel.ActiveWorkbook

$backlog  = $wb.Worksheets.Item("Backlog")
$sprint2  = $wb.Worksheets.Item("Sprint2")
$burndown = $wb.Worksheets.Item("Burndown")

# --- Backlog: register the two new sprint-2 stories -----------------------
$backlog.Range("A11").Value = 2
$backlog.Range("B11").Value = "US01"
$backlog.Range("C11").Value = "Dates before current date"
$backlog.Range("D11").Value = "SK"
$backlog.Range("E11").Value = "Coding"

$backlog.Range("A12").Value = 2
$backlog.Range("B12").Value = "US04"
$backlog.Range("C12").Value = "Marriage before divorce"
$backlog.Range("D12").Value = "SK"
$backlog.Range("E12").Value = "Coding"

# --- Sprint2: add the matching detail rows with estimates ------------------
$sprint2.Range("A2").Value = "US01"
$sprint2.Range("B2").Value = "Dates before current date"
$sprint2.Range("C2").Value = "SK"
$sprint2.Range("D2").Value = "Coding"
$sprint2.Range("E2").Value = 50
$sprint2.Range("F2").Value = 60

$sprint2.Range("A3").Value = "US04"
$sprint2.Range("B3").Value = "Marriage before divorce"
$sprint2.Range("C3").Value = "SK"
$sprint2.Range("D3").Value = "Coding"
$sprint2.Range("E3").Value = 80
$sprint2.Range("F3").Value = 120

# Leave the selection on the newly entered block, like a user would.
$sprint2.Activate()
[void]$sprint2.Range("B2:D3").Select()

# --- Burndown: update this check-in's LOC count -----------------------------
$burndown.Activate()
$burndown.Range("D3").Value = 983
[void]$burndown.Range("F6").Select()

# --- Final view state: Backlog is the active/visible tab -------------------
$backlog.Activate()
